$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the user's e-mail / username in the test data sheet.
$ws.Range("A1").Value = "tallu@gmail.com"
$ws.Range("B1").Value = "tallu"

# Turn the e-mail cell into a mailto hyperlink (adds the Hyperlink style too).
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:tallu@gmail.com")

# Move the active selection to A2, as left by the author after the edit.
$null = $ws.Range("A2").Select()
